# Generate Report for Handoff
#
# This script mirrors the commit "Generate Report for Handoff": the
# Overview sheet gets two additional rows (a markdown handback file and a
# second dependent image) on top of the refreshed original row, and the
# zh-cn / de-de detail sheets gain matching rows that also record the new
# dependency bookkeeping columns (Handoff Reason / Dependency From).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$dateHandoff = "2016-03-21 15:00:32"
$dateDetail  = "2016-03-21 15:00:25"
$dateEpoch   = "0001-01-01 00:00:00"

$pngA   = "9b183233-1b05-4502-81bf-9e3b1c4d4f55.png"
$mdFile = "b3be930c-1501-4fa8-b5f5-2052c36566af.md"
$pngB   = "ee89eddc-06b9-420d-b8d5-d7eb1cdfb5cb.png"

$targetPngRow2 = "2372fbd8306eb0930ec61782350d5cd4bbb36235.png"
$targetPngRow4 = "3eea87ec90f066a9b5651649c07cb33c6d191dfe.png"
$zhXlf         = "b3be930c-1501-4fa8-b5f5-2052c36566af.b943cf799926c3eaff54ac3ccaef8abbb4fd6e0e.zh-cn.xlf"
$deXlf         = "b3be930c-1501-4fa8-b5f5-2052c36566af.b943cf799926c3eaff54ac3ccaef8abbb4fd6e0e.de-de.xlf"

$dependencyFrom = "e2e\b3be930c-1501-4fa8-b5f5-2052c36566af.md"

$ready = "Ready for handoff"

$urlE2E = "https://github.com/OpenLocalizationTest/oltest/blob/5168b1d78b6fec337c64283ccc6043205debd315/e2e/"
$urlZH  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3f8bacb9080d7660e3f09e21bde4bee3a4473216/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$urlDE  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/914a6ca1fd3fa305cbbf9ff41d0e0f5b615d32a5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

function Style-FileNameCell($range) {
    $range.Font.Underline = 2
    $range.Font.Color = 15570276
}

function Style-DateCell($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------

# Refresh the existing handoff row (file renamed from .md to .png, and the
# handoff date moved forward).
$ws1.Range("A2").Value2 = $pngA
$ws1.Range("B2").Value2 = $ready
$ws1.Range("C2").Value2 = $ready
$ws1.Range("D2").Value2 = $dateHandoff
Style-DateCell $ws1.Range("D2")

# New row for the markdown file itself.
$ws1.Range("A3").Value2 = $mdFile
$ws1.Range("B3").Value2 = $ready
$ws1.Range("C3").Value2 = $ready
$ws1.Range("D3").Value2 = $dateHandoff
Style-FileNameCell $ws1.Range("A3")
Style-DateCell $ws1.Range("D3")

# New row for the second dependent image.
$ws1.Range("A4").Value2 = $pngB
$ws1.Range("B4").Value2 = $ready
$ws1.Range("C4").Value2 = $ready
$ws1.Range("D4").Value2 = $dateHandoff
Style-FileNameCell $ws1.Range("A4")
Style-DateCell $ws1.Range("D4")

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $urlE2E + $pngA, "", "", $pngA) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), $urlE2E + $mdFile, "", "", $mdFile) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), $urlE2E + $pngB, "", "", $pngB) | Out-Null

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------

$ws2.Range("A2").Value2 = $pngA
$ws2.Range("B2").Value2 = ".png"
$ws2.Range("C2").Value2 = $ready
$ws2.Range("D2").Value2 = $targetPngRow2
$ws2.Range("E2").Value2 = $dateDetail
$ws2.Range("H2").Value2 = $dateEpoch
$ws2.Range("J2").Value2 = "IsDependency"
$ws2.Range("K2").Value2 = $dependencyFrom
Style-FileNameCell $ws2.Range("D2")
Style-DateCell $ws2.Range("E2")
Style-DateCell $ws2.Range("H2")

$ws2.Range("A3").Value2 = $mdFile
$ws2.Range("B3").Value2 = ".md"
$ws2.Range("C3").Value2 = $ready
$ws2.Range("D3").Value2 = $zhXlf
$ws2.Range("E3").Value2 = $dateDetail
$ws2.Range("H3").Value2 = $dateEpoch
$ws2.Range("J3").Value2 = "Include"
Style-FileNameCell $ws2.Range("A3")
Style-FileNameCell $ws2.Range("D3")
Style-DateCell $ws2.Range("E3")
Style-DateCell $ws2.Range("H3")

$ws2.Range("A4").Value2 = $pngB
$ws2.Range("B4").Value2 = ".png"
$ws2.Range("C4").Value2 = $ready
$ws2.Range("D4").Value2 = $targetPngRow4
$ws2.Range("E4").Value2 = $dateDetail
$ws2.Range("H4").Value2 = $dateEpoch
$ws2.Range("J4").Value2 = "IsDependency"
$ws2.Range("K4").Value2 = $dependencyFrom
Style-FileNameCell $ws2.Range("A4")
Style-DateCell $ws2.Range("E4")
Style-DateCell $ws2.Range("H4")

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $urlE2E + $pngA, "", "", $pngA) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), $urlZH + $targetPngRow2, "", "", $targetPngRow2) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), $urlE2E + $mdFile, "", "", $mdFile) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), $urlZH + $zhXlf, "", "", $zhXlf) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), $urlE2E + $pngB, "", "", $pngB) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D4"), $urlZH + $targetPngRow4, "", "", $targetPngRow4) | Out-Null

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------

$ws3.Range("A2").Value2 = $pngA
$ws3.Range("B2").Value2 = ".png"
$ws3.Range("C2").Value2 = $ready
$ws3.Range("D2").Value2 = $targetPngRow2
$ws3.Range("E2").Value2 = $dateHandoff
$ws3.Range("H2").Value2 = $dateEpoch
$ws3.Range("J2").Value2 = "IsDependency"
$ws3.Range("K2").Value2 = $dependencyFrom
Style-FileNameCell $ws3.Range("D2")
Style-DateCell $ws3.Range("E2")
Style-DateCell $ws3.Range("H2")

$ws3.Range("A3").Value2 = $mdFile
$ws3.Range("B3").Value2 = ".md"
$ws3.Range("C3").Value2 = $ready
$ws3.Range("D3").Value2 = $deXlf
$ws3.Range("E3").Value2 = $dateHandoff
$ws3.Range("H3").Value2 = $dateEpoch
$ws3.Range("J3").Value2 = "Include"
Style-FileNameCell $ws3.Range("A3")
Style-FileNameCell $ws3.Range("D3")
Style-DateCell $ws3.Range("E3")
Style-DateCell $ws3.Range("H3")

$ws3.Range("A4").Value2 = $pngB
$ws3.Range("B4").Value2 = ".png"
$ws3.Range("C4").Value2 = $ready
$ws3.Range("D4").Value2 = $targetPngRow4
$ws3.Range("E4").Value2 = $dateHandoff
$ws3.Range("H4").Value2 = $dateEpoch
$ws3.Range("J4").Value2 = "IsDependency"
$ws3.Range("K4").Value2 = $dependencyFrom
Style-FileNameCell $ws3.Range("A4")
Style-DateCell $ws3.Range("E4")
Style-DateCell $ws3.Range("H4")

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $urlE2E + $pngA, "", "", $pngA) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), $urlDE + $targetPngRow2, "", "", $targetPngRow2) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), $urlE2E + $mdFile, "", "", $mdFile) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), $urlDE + $deXlf, "", "", $deXlf) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), $urlE2E + $pngB, "", "", $pngB) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D4"), $urlDE + $targetPngRow4, "", "", $targetPngRow4) | Out-Null
